# Updated symbol list on Sat Jan  7 17:32:13 UTC 2023 with GitHub Actions
# Applies the latest price/volume snapshot values to the crypto tracking sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "D2"  = "261.07";  "E2"  = "0.92%"
    "D3"  = "27.08";   "E3"  = "0.83%"
    "D4"  = "4.702";   "E4"  = "1.28%"
    "D5"  = "0.06182"; "E5"  = "3.27%"
    "D6"  = "6.685";   "E6"  = "0.66%"
    "D7"  = "0.8510";  "E7"  = "-0.73%"
    "D8"  = "0.9156";  "E8"  = "-0.82%"
    "D9"  = "0.1406";  "E9"  = "1.29%"
    "D10" = "0.04659"; "E10" = "2.73%"
    "D11" = "0.07092"; "E11" = "1.20%"
    "D12" = "0.03154"; "E12" = "3.31%"
    "D13" = "0.09038"; "E13" = "-0.86%"
    "D14" = "0.001535";  "E14" = "0.47%"
    "D15" = "0.0006147"; "E15" = "1.56%"
    "D16" = "0.006032";  "E16" = "-1.60%"
    "D17" = "3.450";   "E17" = "0.24%"
    "D18" = "3.179"
    "D19" = "2.169";   "E19" = "0.68%"
    "E20" = "-0.87%"
    "E21" = "0.89%"
    "D22" = "4.079";   "E22" = "1.27%"
    "D23" = "0.04211"; "E23" = "-0.18%"
    "E24" = "0.01%"
    "E25" = "-5.65%"
    "E26" = "0.09%"
    "D40" = "0.03887"; "E40" = "1.52%"
    "D41" = "0.1110";  "E41" = "-0.18%"
    "D42" = "0.004101"; "E42" = "9.14%"
    "E44" = "-9.69%"
    "D45" = "0.00005156"; "E45" = "0.52%"
    "E46" = "0.10%"
    "E48" = "-10.81%"
    "E49" = "0.10%"
    "E50" = "0.10%"
}

foreach ($addr in $updates.Keys) {
    # Leading apostrophe forces Excel to store these numeric-looking
    # strings (prices / percentages) as text, matching the source sheet
    # where column D/E are plain text cells, not numbers.
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
